# Added ifoCAST full series evaluation: the quarter-over-quarter error series
# (row 2 = Q0 .. row 24) is recomputed with one more forecast quarter included.
# Effect on the sheet: for each data row the existing error values shift one
# column to the left (B<-C, C<-D, ... ) to make room for the newly evaluated
# quarter, a fresh error value is appended at the end of the rows that still
# have a full 10-quarter window (rows 2-14), and the rows that were already
# shrinking (rows 15-24, the triangular "tail" of the matched-error table)
# simply lose their last column instead of gaining a new value.
$rowData = @{}
$rowData[2] = @(9.141635797313464, -8.409018285787846, -0.6996216861316987, 0.8832909747552129, -2.082996017594043, 0.18103513707042, -0.2701190178563186, -0.4104830813702928, 0.1157559294919248, -0.02607301724734923)
$rowData[3] = @(-8.723943454208817, -1.01454685455267, 0.5683658063342414, -2.397921186015015, -0.1338900313505515, -0.5850441862772902, -0.7254082497912643, -0.1991692389290468, -0.3409981856683208, 0.3019764357337431)
$rowData[4] = @(-2.0709218147352, -0.4880091538482882, -3.454296146197544, -1.190264991533081, -1.64141914645982, -1.781783209973794, -1.255544199111576, -1.39737314585085, -0.7543985244487865, -0.6299858181192041)
$rowData[5] = @(0.3826372221226423, -2.583649770226613, -0.3196186155621505, -0.7707727704888893, -0.9111368340028634, -0.3848978231406458, -0.5267267698799198, 0.1162478515221441, 0.2406605578517265, -0.1569696481858908)
$rowData[6] = @(-2.594300986556497, -0.3302698318920341, -0.7814239868187727, -0.9217880503327469, -0.3955490394705293, -0.5373779862098034, 0.1055966351922606, 0.230009341521843, -0.1676208645157742, -0.24378957821662)
$rowData[7] = @(-0.4388843641081749, -0.8900385190349136, -1.030402582548888, -0.5041635716866701, -0.6459925184259441, -0.003017897023880223, 0.1213948093057022, -0.2762353967319151, -0.3524041104327608, -0.07608205230983922)
$rowData[8] = @(-0.6819367980671746, -0.8223008615811487, -0.2960618507189311, -0.4378907974582051, 0.2050838239438588, 0.3294965302734412, -0.06813367576417605, -0.1443023894650218, 0.1320196686578998, -0.0392001536538242)
$rowData[9] = @(-0.7058026912073482, -0.1795636803451306, -0.3213926270844047, 0.3215819943176592, 0.4459947006472416, 0.04836449460962439, -0.02780421909122137, 0.2485178390317002, 0.07729801671997623, 0.2304454222287174)
$rowData[10] = @(-0.1340702663479956, -0.2758992130872696, 0.3670754083147943, 0.4914881146443768, 0.09385790860675949, 0.01768919490591373, 0.2940112530288354, 0.1227914307171113, 0.2759388362258526, -0.001489938197266189)
$rowData[11] = @(-0.293453334447234, 0.3495212869548299, 0.4739339932844123, 0.07630378724679503, 0.0001350735459492769, 0.2764571316688709, 0.1052373093571469, 0.2583847148658881, -0.01904405955723064, -0.182031752916177)
$rowData[12] = @(0.4376286059058361, 0.5620413122354185, 0.1644111061978012, 0.08824239249695551, 0.3645644506198771, 0.1933446283081531, 0.3464920338168943, 0.06906325939377558, -0.09392443396517081, -0.2180070093596886)
$rowData[13] = @(0.7406159457232021, 0.3429857396855849, 0.2668170259847391, 0.5431390841076607, 0.3719192617959367, 0.525066667304678, 0.2476378928815592, 0.0846501995226128, -0.03943237587190501, 0.4767206611340558)
$rowData[14] = @(1.255012967438235, 1.178844253737389, 1.455166311860311, 1.283946489548587, 1.437093895057328, 1.159665120634209, 0.9966774272752628, 0.8725948518807449, 1.388747888886706, 1.178844253737389)
$rowData[15] = @(0.2348700177716323, 0.5111920758945538, 0.3399722535828299, 0.4931196590915711, 0.2156908846684524, 0.05270319130950599, -0.07137938408501182, 0.444773652920949, 0.2348700177716323)
$rowData[16] = @(0.5151599734076631, 0.343940151095939, 0.4970875566046802, 0.2196587821815615, 0.0566710888226151, -0.06741148657190271, 0.4487415504340581, 0.2388379152847414)
$rowData[17] = @(0.4795802412661804, 0.6327276467749217, 0.3552988723518029, 0.1923111789928565, 0.06822860359833866, 0.5843816406042994, 0.3744780054549828)
$rowData[18] = @(0.3919214649192569, 0.1144926904961382, -0.04849500286280822, -0.172577578257326, 0.3435754587486348, 0.1336718235993181)
$rowData[19] = @(0.0691614752440418, -0.09382621811490459, -0.2179087935094224, 0.2982442434965384, 0.08834060834722172)
$rowData[20] = @(-0.1606876400509585, -0.2847702154454763, 0.2313828215604846, 0.02147918641116785)
$rowData[21] = @(-0.3143564178021929, 0.201796619203768, -0.00810701594554874)
$rowData[22] = @(0.1836459624741271, -0.02625767267518964)
$rowData[23] = @(-0.04428949692388896)
$rowData[24] = @()

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 24; $r++) {
    $vals = $rowData[$r]
    $n = $vals.Count
    # Write new values starting at column B (col 2)
    for ($i = 0; $i -lt $n; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $vals[$i]
    }
    # Clear any remaining cells from column (2+n) through K (col 11)
    for ($c = 2 + $n; $c -le 11; $c++) {
        $ws.Cells.Item($r, $c).ClearContents()
    }
}
